$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.442.53"
$ws.Range("E2").Value = "  +1.11%  "

$ws.Range("D3").Value = "2.435.44"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.45%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.510"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.17%  "

$ws.Range("E9").Value = "  +9.57%  "

$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.327"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.10%  "

$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.26%  "

$ws.Range("D13").Value = "68.331.29"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("E14").Value = "  +4.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "336.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("E19").Value = "  +2.30%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("B21").Value = "SuiNetwork"
$ws.Range("C21").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.37%  "

$ws.Range("E23").Value = "  +2.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.51%  "

$ws.Range("D25").Value = "0.0₃0811"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "424.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.45%  "

$ws.Range("E29").Value = "  +3.12%  "

$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.65%  "

$ws.Range("E35").Value = "  -1.76%  "

$ws.Range("E36").Value = "  -0.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.14%  "

$ws.Range("E38").Value = "  +3.52%  "

$ws.Range("E39").Value = "  -0.14%  "

$ws.Range("E40").Value = "  +1.66%  "

$ws.Range("E41").Value = "  +2.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "129.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0715"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.479"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.560"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0917"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.40%  "

$ws.Range("E47").Value = "  +1.17%  "

$ws.Range("E48").Value = "  -1.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").Value = "0.0₆0203"
$ws.Range("E51").Value = "  +4.62%  "

